$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("B8").Value = "-"
